# Ticket 57 - Add "MultiLevel2" and "Grid" worksheets with sample data used
# by the new nested-loop / OFFSET-based formula regression tests, and make
# "Grid" the active sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the "MultiLevel2" worksheet (after the last existing sheet)
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws7 = $wb.Worksheets.Add($null, $lastSheet)
$ws7.Name = "MultiLevel2"

# --- Write cell values first, in the same order the template was
#     originally authored, so shared-string indices line up. ---
$ws7.Range("A1").Value = 'Department'
$ws7.Range("A2").Value = 'Installation'
$ws7.Range("B1").Value = 'Job Cost'
$ws7.Range("C1").Value = 'Materials Cost'
$ws7.Range("D1").Value = 'Total Cost'
$ws7.Range("A3").Value = '<jt:forEach items="${workOrders}" var="workOrder" groupBy="department;location">${workOrder.obj.department} ${workOrder.obj.location}'
$ws7.Range("A4").Value = '<jt:forEach items="${workOrder.items}" var="wo" groupBy="installation"><jt:forEach items="${wo.items}" var="detail">${detail.installation}'
$ws7.Range("B4").Value = '${detail.jobAmt}'
$ws7.Range("C4").Value = '${detail.matAmt}'
$ws7.Range("E4").Value = '</jt:forEach>'
$ws7.Range("A5").Value = '${wo.obj.installation} Total:'
$ws7.Range("B5").Value = '$[SUM(B4)]'
$ws7.Range("C5").Value = '$[SUM(C4)]'
$ws7.Range("D5").Value = '$[SUM(D4)]'
$ws7.Range("A7").Value = '${workOrder.obj.department} ${workOrder.obj.location} Total:'
$ws7.Range("B7").Value = '$[SUM(B5)]'
$ws7.Range("C7").Value = '$[SUM(C5)]'
$ws7.Range("D7").Value = '$[SUM(D5)]'
$ws7.Range("A9").Value = 'Grand Total:'
$ws7.Range("B9").Value = '$[SUM(B7)]'
$ws7.Range("C9").Value = '$[SUM(C7)]'
$ws7.Range("D9").Value = '$[SUM(D7)]'
$ws7.Range("D4").Value = '${detail.totAmt}'
$ws7.Range("E6").Value = '</jt:forEach>'
$ws7.Range("E7").Value = '</jt:forEach>'

# --- Column widths ---
$ws7.Columns.Item(1).ColumnWidth = 35.43
$ws7.Columns.Item(2).ColumnWidth = 11.29
$ws7.Columns.Item(3).ColumnWidth = 13.71
$ws7.Columns.Item(4).ColumnWidth = 12.57

# --- Row heights for the thin spacer rows between groups ---
$ws7.Rows.Item(6).RowHeight = 4.5
$ws7.Rows.Item(8).RowHeight = 4.5

# --- Number formats (currency) ---
$ws7.Range("B1:D2").NumberFormat = '"$"#,##0.00'
$ws7.Range("B4:E4").NumberFormat = '"$"#,##0.00'
$ws7.Range("B5:D5").NumberFormat = '"$"#,##0.00'
$ws7.Range("B7:E7").NumberFormat = '"$"#,##0.00'
$ws7.Range("B9:D9").NumberFormat = '"$"#,##0.00'

# --- Fonts / alignment ---
$ws7.Range("A1:D2").Font.Bold = $true
$ws7.Range("A3").Font.Bold = $true
$ws7.Range("A7").Font.Bold = $true
$ws7.Range("A9").Font.Bold = $true

$ws7.Range("A7").HorizontalAlignment = -4152
$ws7.Range("A9").HorizontalAlignment = -4152
$ws7.Range("A5").HorizontalAlignment = -4152

$ws7.PageSetup.Orientation = 1

# ---------------------------------------------------------------------
# 2. Add the "Grid" worksheet (after "MultiLevel2")
# ---------------------------------------------------------------------
$ws8 = $wb.Worksheets.Add($null, $ws7)
$ws8.Name = "Grid"

$ws8.Range("C2").Value = '$[SUM(B2)]'
$ws8.Range("C3").Value = '$[SUM(C2)]'
$ws8.Range("A1").Value = 'Region'
$ws8.Range("B1").Value = '<jt:forEach items="${dates}" var="dateString" copyRight="true" indexVar="i">${dateString}</jt:forEach>'
$ws8.Range("A2").Value = '<jt:forEach items="${regions}" var="region">${region.name}'
$ws8.Range("B2").Value = '<jt:forEach items="${region.salesFigures}" var="sales" copyRight="true">${sales}</jt:forEach>'
$ws8.Range("B3").Value = '<jt:forEach items="${dates}" var="date" copyRight="true" indexVar="i"><jt:formula text="SUM(OFFSET(B2, 0, ${i}, ${regions.size()}, 1))"/></jt:forEach>'
$ws8.Range("C1").Value = 'Totals'
$ws8.Range("A3").Value = 'Totals'
$ws8.Range("D2").Value = '</jt:forEach>'

$ws8.Columns.Item(1).ColumnWidth = 26

$ws8.Range("A1:C1").Font.Bold = $true
$ws8.Range("C2").Font.Bold = $true
$ws8.Range("A3:C3").Font.Bold = $true

$ws8.Range("R15").NumberFormat = '"$"#,##0.00'

# ---------------------------------------------------------------------
# 3. Make "Grid" the active sheet/tab
# ---------------------------------------------------------------------
$ws8.Activate()
$ws8.Range("A1").Select()
